$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write values column-by-column (B then C then D) to preserve the same
# shared-string insertion order as the original workbook. Purely numeric
# looking values get a leading apostrophe (quote-prefix) so Excel stores
# them as text rather than converting them to numbers - matching the
# shared-string <t> entries from the source workbook.
$ws.Range("B2").Value = "'0.17"
$ws.Range("B3").Value = "'-0.01"
$ws.Range("B4").Value = "'-0.09"

$ws.Range("C2").Value = "44.29***"
$ws.Range("C3").Value = "2.21***"
$ws.Range("C4").Value = "'0.98"

$ws.Range("D2").Value = "'-0.89"
$ws.Range("D3").Value = "0.46***"
$ws.Range("D4").Value = "0.82*"
